$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.708.94"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.301.00"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.28"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.75"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.56"
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.15"
$ws.Range("E11").Value = "  +5.25%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.657.57"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.290.50"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.661.52"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.25"
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.68"
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.29"
$ws.Range("E23").Value = "  +8.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.90"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.41"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.33"
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("E28").Value = "  +15.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.50"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.06"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.93"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.72"
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("E35").Value = "  -5.15%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0695"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.100"
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.108"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.69"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.29"
$ws.Range("E42").Value = "  +13.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.950.39"
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.30"
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("E46").Value = "  +3.87%  "
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.526.76"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.35"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.49"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("E51").Value = "  +1.17%  "
